# Fruta / hortaliza, semanal
# The weekly refresh reshuffles which historical record lives on which row:
# each target row picks up the Fecha/Calidad/Volumen/Precio.../Unidad/Precio $/Kg/
# Kg-unidad values that used to belong to a different row (Mercado/Región/
# Producto/Categoría/Variedad stay put since they describe the same series).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as a block, keyed by column letter.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "S", "T")

# Snapshot every relevant cell in rows 2-15 (row 9 and row 16 are untouched)
# before any writes happen, since several rows trade values with each other.
$snapshot = @{}
foreach ($r in 2..15) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Destination row -> source row (source row's old values become the
# destination row's new values).
$mapping = @{
    2  = 15
    3  = 12
    4  = 7
    5  = 14
    6  = 5
    7  = 4
    8  = 6
    10 = 3
    11 = 2
    12 = 13
    13 = 8
    14 = 11
    15 = 10
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
